$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: remove the existing "_GoBack" bookmark that sits right after
# the "Пользовательский " run (it gets relocated further down the doc).
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# Change 2: "Запоминание последнего населенного пункта для
# автоматической переадресации" becomes three runs:
#   "Запоминание последнего" | " просмотренного" | (bookmark _GoBack) |
#   " населенного пункта для автоматической переадресации"
# We split the run in two steps using temporary bookmarks (bookmarks
# force a run boundary at the collapsed insertion point), insert the
# new text, re-add the real "_GoBack" bookmark at the right spot, then
# clean up the scaffolding bookmarks (removing a bookmark does not
# re-merge the runs it separated).
# ---------------------------------------------------------------------

# 2a. Mark the point right after "...последнего" (before " населенного...")
$rng = $d.Content
$rng.Find.Execute("Запоминание последнего") | Out-Null
$pointA = $rng.Duplicate
$pointA.Collapse(0)
$d.Bookmarks.Add("TEMP_SPLIT_A", $pointA) | Out-Null

# 2b. Insert the new text " просмотренного" right before that marker.
$bmA = $d.Bookmarks("TEMP_SPLIT_A")
$bmA.Range.InsertBefore(" просмотренного")

# 2c. Mark the point right after "...последнего" again (now right before
#     the newly inserted " просмотренного"), splitting that run too.
$rng2 = $d.Content
$rng2.Find.Execute("Запоминание последнего") | Out-Null
$pointB = $rng2.Duplicate
$pointB.Collapse(0)
$d.Bookmarks.Add("TEMP_SPLIT_B", $pointB) | Out-Null

# 2d. Re-create the real "_GoBack" bookmark exactly where TEMP_SPLIT_A is.
$bmA2 = $d.Bookmarks("TEMP_SPLIT_A")
$d.Bookmarks.Add("_GoBack", $bmA2.Range) | Out-Null

# 2e. Drop the scaffolding bookmarks; the run split they created stays.
$d.Bookmarks("TEMP_SPLIT_A").Delete()
$d.Bookmarks("TEMP_SPLIT_B").Delete()

# ---------------------------------------------------------------------
# Change 3: the run sequence " с" + "обрав статистику" + "," +
# " можно " + "с" + "делать прогнозы на стихийные бедствия. Можно
# продавать МЧС" + " или " collapses into a single run. A no-op
# replace (same text in, same text out) on an interior substring makes
# the engine re-normalize/merge the contiguous same-formatted runs
# around it, without touching the separately formatted
# "Экономическая выгода:" run before it.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("обрав статистику", $true, $false, $false, $false, $false, $true, 1, $false, "обрав статистику", 2) | Out-Null

Write-Host "edit complete"
